# Revert "UPF IP overlap" (commit aa27073a64bc6f81777792c741985d00ed590c2d)
# on the single content slide. The original edit had: shrunk the first
# "Rectangulo: esquinas redondeadas" grouping box, moved the "2" icon/label
# pair over to the UPF box, and added a brand-new dashed rounded rectangle
# around UPF. Reverting restores the original wider box, puts the "2"
# icon/label back at its original spot, and removes the extra rectangle.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Shape id=3 "Rectangulo: esquinas redondeadas 2" - restore original (wider) width
$rect2 = Get-ShapeById $s 3
$rect2.Width = 366.7755905511811

# Shape id=35 "Picture 2" - move back to its original spot next to UPF
$pic35 = Get-ShapeById $s 35
$pic35.Left = 182.89708711417322
$pic35.Top = 205.43244174488188

# Shape id=36 "CuadroTexto 35" (the "2" label) - move back alongside the picture
$txt36 = Get-ShapeById $s 36
$txt36.Left = 184.7732315464567
$txt36.Top = 205.43244174488188

# Shape id=47 "Rectangulo: esquinas redondeadas 46" - this shape was newly
# added by the "UPF IP overlap" commit; delete it to complete the revert.
$rect47 = Get-ShapeById $s 47
$rect47.Delete()
